$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.906.91'
$ws.Range("E2").Value = '  -1.67%  '

$ws.Range("D3").Value = '2.906.00'
$ws.Range("E3").Value = '  -2.50%  '

$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '580.29'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.48%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.74'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.55%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.12%  '

$ws.Range("E8").Value = '  -0.17%  '

$ws.Range("D9").Value = '2.903.93'
$ws.Range("E9").Value = '  -2.52%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.73'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -7.28%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.153'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.34%  '

$ws.Range("E12").Value = '  -3.04%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000240'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.01%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.67'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.01%  '

$ws.Range("E15").Value = '  -0.81%  '

$ws.Range("D16").Value = '3.387.00'
$ws.Range("E16").Value = '  -2.27%  '

$ws.Range("D17").Value = '61.868.02'
$ws.Range("E17").Value = '  -1.46%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.68'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.76%  '

$ws.Range("D19").Value = '2.918.60'
$ws.Range("E19").Value = '  -1.82%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '435.96'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.86%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.35'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.09%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.660'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.97%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.96'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.06%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '80.00'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.28%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.02'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.85%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.24'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -9.56%  '

$ws.Range("E27").Value = '  -0.03%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.06'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.41%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000111'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +15.63%  '

$ws.Range("E30").Value = '  -1.68%  '

$ws.Range("E31").Value = '  -2.36%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.11'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.28%  '

$ws.Range("E33").Value = '  -1.95%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.998'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.23%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '25.79'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.70%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.965'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.55%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.09'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.43%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.50'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.55%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '49.12'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.03%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.98'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.01%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.33'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.77%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.116'
$ws.Range("D42").Style = "Normal"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.271'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.93%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '38.32'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.97%  '

$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '2.687.43'
$ws.Range("E45").Value = '  -2.22%  '

$ws.Range("B46").Value = 'Monero'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '134.33'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.80%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0336'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.53%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '343.86'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.64%  '

$ws.Range("E50").Value = '  -1.92%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '22.00'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.91%  '
